$d = $word.ActiveDocument

# --- helper: split a run in two (same formatting) at a given prefix boundary,
# without altering any visible text, by toggling Bold on/off across the
# prefix sub-range. Word always breaks run boundaries at the edge of any
# range whose formatting is (re)applied, even when the applied value is a
# no-op, which is exactly the "retyped/reflowed" run-split pattern produced
# by a PDF-text-extraction edit.
function Split-Run($prefixText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($prefixText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $found) {
        throw "Split-Run: text not found: $prefixText"
    }
    $rng.Bold = 1
    $rng.Bold = 0
}

# 1) Body paragraph: "...concerned that policymak" | "ers are proposing..."
Split-Run("I am a former employee of Illinois Tool Works and have received Medicare Advantage coverage since retiring. As a senior who relies on my social security and pension checks to make ends meet, I am really concerned that policymak")

# 2) Body paragraph: "...benefits. I also" | " enjoy reading the booklets..."
Split-Run("I am very pleased with my WellMed Medicare Advantage plan. My health care providers and insurance company make sure that I have a great understanding of my benefits. I also")

# 3) Body paragraph: "...benefits should re" | "main stable for all seniors..."
Split-Run("I believe that Medicare Advantage benefits should re")

# 4) Body paragraph: "...funding for Medicare Advantage" | " and that they really..."
Split-Run("I hope policymakers continue to protect funding for Medicare Advantage")

# --- bookmark "_GoBack" spanning from the very start of the document to
# just after the final visible character (the page-number "1"), matching
# the position Word leaves it at after the last edit.
$startRng = $d.Range(0, 0)
$lastPara = $d.Paragraphs.Last
$endPos = $lastPara.Range.End - 1
$bmRng = $d.Range($startRng.Start, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# --- sectPr: restore header/footer distance, gutter and column spacing
# defaults that a full Word round-trip always re-emits.
$sec = $d.Sections.First
$sec.PageSetup.HeaderDistance = 36
$sec.PageSetup.FooterDistance = 36
$sec.PageSetup.Gutter = 0
$sec.PageSetup.TextColumns.Spacing = 36

Write-Output "done"
